# Ninth commit with signup and login page test fixed
# Adds a new test case (TC_07 - Duplicate Email Registration Check) to the
# Nesto test-cases worksheet, appended as rows 27-35.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 27: new test case header row (Test Case ID, Description, first step)
$ws.Range("A27").Value = "TC_07"
$ws.Range("B27").Value = "Duplicate Email Registration Check"
$ws.Range("C27").Value = '1.Open URL "http://localhost:8080/signup"'

# Rows 28-32: individual test steps
$ws.Range("C28").Value = '2.Type "Duplicate Test" into "//input[@placeholder=''Enter Full Name'']"'
$ws.Range("C29").Value = '3.Type "faizal@nesto.com" into "//input[@placeholder=''admin@nesto.com'']"'
$ws.Range("C30").Value = '4.Type "anyPassword123" into "//input[@placeholder=''Create Password'']"'
$ws.Range("C31").Value = '5.Type "32" into "//input[@placeholder=''Your Age'']"'
$ws.Range("C32").Value = '6.Type "9947110008" into "//input[@placeholder=''Your Mobile Number'']"'

# Row 33: reuse of existing step text ("7.Click on the Sign Up button ...")
$ws.Range("C33").Value = '7.Click on the "Sign Up" button "//button[@type=''submit'']"'

# Row 35 is populated before row 34 to mirror the original authoring order
# of the shared-strings table (string "9.Verify text..." was interned
# ahead of "8.Verify URL contains...").
$ws.Range("C35").Value = '9.Verify text "Email already registered. Please login." at "//div[@class=''error-alert'']"'

# Row 34: verification step
$ws.Range("C34").Value = '8.Verify URL contains "register"'

# Update the visible selection/scroll position to match the authored view
$ws.Application.ActiveWindow.ScrollRow = 11
$ws.Range("C34").Select()
